$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 126
$ws.Range("A126").Value = "insert_027"
$ws.Range("B126").Value = "y"
$ws.Range("C126").Value = "布尔型字段插入小数"
$ws.Range("D126").Value = "insert"
$ws.Range("F126").Value = "schema1"
$ws.Range("G126").Value = "insert_value22"
$ws.Range("H126").Value = "4"
$ws.Range("I126").Value = "select * from `$schema1"
$ws.Range("J126").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/insert/expectedresult/insert_027.csv"
$ws.Range("O126").Value = "csv_containsAll"

# Row 127
$ws.Range("A127").Value = "array_05"
$ws.Range("B127").Value = "y"
$ws.Range("C127").Value = "布尔型数组插入元素为小数"
$ws.Range("D127").Value = "ComplexDataType"
$ws.Range("E127").Value = "Array"
$ws.Range("F127").Value = "array10"
$ws.Range("G127").Value = "array10_value58"
$ws.Range("H127").Value = "3"
$ws.Range("I127").Value = "select in_use from `$array10"
$ws.Range("J127").Value = "src/test/resources/io.dingodb.test/testdata/cases/dml/insert/expectedresult/complexdatatype/array/array_005.csv"
$ws.Range("O127").Value = "csv_containsAll"

$ws.Range("A127").Select()
